$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New headers for the two extra metric columns ---
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy the header formatting (bold, centered, bordered) from F1 onto G1:H1
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Updated regression metrics (B: MSE, C: R2, D: MAE) plus new
#     Elapsed Time (G) / CPU (H) columns for every data row ---
$metrics = @(
    @{ Row = 2;  B = 0.4972636062298078;  C = 0.9900984256003985;  D = 0.5737939481677846 },
    @{ Row = 3;  B = 0.2376395757976146;  C = 0.9953553285148486;  D = 0.3829637447440828 },
    @{ Row = 4;  B = 0.3440597266435089;  C = 0.9933810559248661;  D = 0.4635476650442505 },
    @{ Row = 5;  B = 0.4134410054541373;  C = 0.9918475094236497;  D = 0.5003202100897108 },
    @{ Row = 6;  B = 0.6666393828293923;  C = 0.9804190550216592;  D = 0.6111898016761881 },
    @{ Row = 7;  B = 0.2476002055963841;  C = 0.9965780981173442;  D = 0.3764936609166029 },
    @{ Row = 8;  B = 0.1533455515139108;  C = 0.9984195000978843;  D = 0.3352655043891136 },
    @{ Row = 9;  B = 0.4749149041100557;  C = 0.9971710674202325;  D = 0.5754557717508432 },
    @{ Row = 10; B = 0.1075861073149535;  C = 0.9980346566053806;  D = 0.2435166343840025 },
    @{ Row = 11; B = 0.2540151239672845;  C = 0.9812324372710919;  D = 0.3884996060380089 },
    @{ Row = 12; B = 0.05043507891103303; C = 0.9985167903739312;  D = 0.1647042703662659 },
    @{ Row = 13; B = 0.09783325630578199; C = 0.9990720077201612;  D = 0.2199165099218983 },
    @{ Row = 14; B = 0.08744915486079569; C = 0.9988159924238357;  D = 0.2401551473425225 }
)

$elapsedTime = 1.669922641383406
$cpu = 0.97

foreach ($m in $metrics) {
    $r = $m.Row
    $ws.Range("B$r").Value = $m.B
    $ws.Range("C$r").Value = $m.C
    $ws.Range("D$r").Value = $m.D
    $ws.Range("G$r").Value = $elapsedTime
    $ws.Range("H$r").Value = $cpu
}

Write-Host "done"
